$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected (no known plaintext password), so unlock it,
# apply the edits, then restore protection.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer note.
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.
Model holdings provided as of 2021-06-09 for illustrative purposes only and are subject to change."
# Undo the implicit row auto-height bump caused by writing a multi-line
# value so row 11 keeps its original (default / non-custom) height.
$ws.Rows(11).EntireRow.AutoFit()

# Refresh the Weight / Percent Change figures for each holding.
$ws.Range("D2").Value = 0.5411981084181636
$ws.Range("E2").Value = -0.003186404673393484

$ws.Range("D3").Value = 0.2503108778370013
$ws.Range("E3").Value = 0.0002885586495455783

$ws.Range("D4").Value = 0.04993701561414619
$ws.Range("E4").Value = -0.004474829086388965

$ws.Range("D5").Value = 0.09932441087450769
$ws.Range("E5").Value = -0.007726570311789915

$ws.Range("D6").Value = 0.02908285442189512
$ws.Range("E6").Value = -0.008511408909815166

$ws.Range("D7").Value = 0.03014673283428608
$ws.Range("E7").Value = -0.007412398921832986

$ws.Range("E8").Value = -0.003114139143483907

$ws.Protect()
